$d = $word.ActiveDocument

function Replace-Exact([string]$old, [string]$new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "MISS: $old"
    }
    return $ok
}

# ---- Edit 1: remove "3.2.2 аппаратным" section; renumber "3.2.3 программным" heading -> "3.2.2" ----
# Paragraph 87 carries the <w:lastRenderedPageBreak/>, so keep that paragraph as the
# container for the renumbered heading and delete the two paragraphs that follow it
# (the COM-port support paragraph, and the old "3.2.3" heading paragraph).
$p87 = $d.Paragraphs(87)
$headingRange = $d.Range($p87.Range.Start, $p87.Range.End - 1)
$headingRange.Text = "3.2.2. Требования к программным интерфейсам"

$delRange = $d.Range($d.Paragraphs(88).Range.Start, $d.Paragraphs(89).Range.End)
$delRange.Delete()

Write-Output "edit1 done"
